# DD - Dizionario dei Dati.xlsx edit
#
# Changes applied (per commit message "Apportate modifiche al Dizionario per
# l'attributo DD_FAT: Cancellato Orario inizio consegna/Orario fine consegna"):
#   1. Remove the two rows in the DD_Fat section that describe
#      "orario di inizio disponibilità" / "orario di fine disponibilità"
#      (rows 28 and 29), shifting every following row up by two.
#   2. Update the length/format description for "via e numero civico" in the
#      DD_Ord section (moved from B39 to B37 after the shift) from
#      "almeno 2 caratteri massimo 23 totali" to
#      "almeno 1 caratteri massimo 20 totali".
#   3. Update the sheet view so it is scrolled back to the top (no
#      topLeftCell override) and the active selection is the new last
#      attribute length cell, B31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the "orario di inizio/fine disponibilità" rows from the DD_Fat
#    block; everything below shifts up automatically.
$ws.Rows("28:29").Delete()

# 2. Fix the "via e numero civico" length description (now at row 37).
$ws.Range("B37").Value = "almeno 1 caratteri massimo 20 totali"

# 3. Reset the view: select B31 (resets scroll position / topLeftCell) and
#    leave it as the active selection, matching the saved workbook state.
$ws.Range("B31").Select()
